$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Mark the "teclas acceso rapido" task as "en proceso"
$ws.Range("C2").Value = "en proceso"

# Update selection to A3 (matches the committed sheetView selection)
$ws.Activate()
$ws.Range("A3").Select()
